$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# 1) Reassign the Module column (H) for the existing rows 3-8 from
#    "MA-20" to "ICT-431" (both strings already exist in the workbook).
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = "ICT-431"
$ws.Range("H4").Value = "ICT-431"
$ws.Range("H5").Value = "ICT-431"
$ws.Range("H6").Value = "ICT-431"
$ws.Range("H7").Value = "ICT-431"
$ws.Range("H8").Value = "ICT-431"

# ---------------------------------------------------------------------------
# 2) Fill in the new log rows (13, 14, 16-20). Writes are ordered so that
#    newly-introduced shared strings are appended in the same sequence as
#    the target workbook (continuer le code, j'ai continuer..., Le prof a
#    parlé..., nous avons fait..., faire une maquette HUD, presenter la
#    maquette HUD, j'ai presenter HUD..., faire une maquette HUD 2.0).
# ---------------------------------------------------------------------------

# Row 13
$ws.Range("A13").Value = 44258
$ws.Range("B13").Value = 0.5625
$ws.Range("C13").Value = 0.59375
$ws.Range("E13").Value = "Pratique"
$ws.Range("F13").Value = "continuer le code"
$ws.Range("G13").Value = "Travail"
$ws.Range("H13").Value = "MA-20"
$ws.Range("I13").Value = "j'ai continuer le code de la bataille navale, j'ai commencer le jeu"
$ws.Range("J13").Value = "NON"

# Row 14
$ws.Range("A14").Value = 44258
$ws.Range("B14").Value = 0.59722222222222221
$ws.Range("C14").Value = 0.62847222222222221
$ws.Range("E14").Value = "Pratique"
$ws.Range("F14").Value = "continuer le code"
$ws.Range("G14").Value = "Travail"
$ws.Range("H14").Value = "MA-20"
$ws.Range("I14").Value = "j'ai continuer le code de la bataille navale, j'ai commencer le jeu"
$ws.Range("J14").Value = "NON"

# Row 16
$ws.Range("A16").Value = 44259
$ws.Range("B16").Value = 0.33333333333333331
$ws.Range("C16").Value = 0.35069444444444442
$ws.Range("E16").Value = "Théorie"
$ws.Range("F16").Value = "ecouter le prof"
$ws.Range("G16").Value = "Travail"
$ws.Range("H16").Value = "ICT-431"
$ws.Range("I16").Value = "Le prof a parlé sur Les Maquettes"
$ws.Range("J16").Value = "OUI"

# Row 17
$ws.Range("A17").Value = 44259
$ws.Range("B17").Value = 0.35416666666666669
$ws.Range("C17").Value = 0.36458333333333331
$ws.Range("E17").Value = "Pratique"
$ws.Range("I17").Value = "nous avons fait une maquette par groupe de quatre"
$ws.Range("F17").Value = "faire une maquette HUD"
$ws.Range("G17").Value = "Travail"
$ws.Range("H17").Value = "ICT-431"
$ws.Range("J17").Value = "NON"

# Row 18
$ws.Range("A18").Value = 44259
$ws.Range("B18").Value = 0.36805555555555558
$ws.Range("C18").Value = 0.39930555555555558
$ws.Range("E18").Value = "Pratique"
$ws.Range("F18").Value = "faire une maquette HUD"
$ws.Range("G18").Value = "Travail"
$ws.Range("H18").Value = "ICT-431"
$ws.Range("I18").Value = "nous avons fait une maquette par groupe de quatre"
$ws.Range("J18").Value = "OUI"

# Row 19
$ws.Range("A19").Value = 44259
$ws.Range("B19").Value = 0.40625
$ws.Range("C19").Value = 0.44444444444444442
$ws.Range("E19").Value = "Pratique"
$ws.Range("F19").Value = "presenter la maquette HUD"
$ws.Range("G19").Value = "Travail"
$ws.Range("H19").Value = "ICT-431"
$ws.Range("I19").Value = "j'ai presenter HUD et ecouter les autres qui on presenter le leur"
$ws.Range("J19").Value = "OUI"

# Row 20
$ws.Range("A20").Value = 44259
$ws.Range("B20").Value = 0.4458333333333333
$ws.Range("C20").Value = 0.45833333333333331
$ws.Range("E20").Value = "Pratique"
$ws.Range("F20").Value = "faire une maquette HUD 2.0"
$ws.Range("G20").Value = "Travail"
$ws.Range("H20").Value = "ICT-431"

# ---------------------------------------------------------------------------
# 3) Column widths (F and I grew to fit the new, longer text).
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 30.830729166666668
$ws.Columns.Item(9).ColumnWidth = 55.166666666666664

# ---------------------------------------------------------------------------
# 4) Move the selection cursor to I20, matching the saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("I20").Select()
